$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 50090.95
$ws.Range("I33").Value = 58915.65
$ws.Range("J33").Value = 84.333336
$ws.Range("K33").Value = 58915.65
$ws.Range("L33").Value = 84.333336
$ws.Range("M33").Value = -58686.65
$ws.Range("N33").Value = -542.333336

$ws.Range("H55").Value = 263.27274
$ws.Range("I55").Value = 221.77777
$ws.Range("J55").Value = 450
$ws.Range("K55").Value = 221.77777
$ws.Range("L55").Value = 450
$ws.Range("M55").Value = -7.777770000000004
$ws.Range("N55").Value = -878

$ws.Range("H126").Value = 38000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 38000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 38000
$ws.Range("N126").Value = -47880

$ws.Range("H129").Value = 1458.6666
$ws.Range("I129").Value = 298
$ws.Range("J129").Value = 1845.5555
$ws.Range("K129").Value = 894
$ws.Range("L129").Value = 5536.666499999999
$ws.Range("M129").Value = 4106
$ws.Range("N129").Value = -15536.6665

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value = 1638.2433
$ws.Range("I132").Value = 1733.909
$ws.Range("J132").Value = 849
$ws.Range("K132").Value = 5201.727000000001
$ws.Range("L132").Value = 2547
$ws.Range("M132").Value = -2671.727000000001
$ws.Range("N132").Value = -7607

$ws.Range("H136").Value = 53777.777
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 53777.777
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 53777.777
$ws.Range("N136").Value = -63977.777

$ws.Range("H137").Value = 1904.7142
$ws.Range("I137").Value = 1629
$ws.Range("J137").Value = 2304.5
$ws.Range("K137").Value = 4887
$ws.Range("L137").Value = 6913.5
$ws.Range("M137").Value = -2337
$ws.Range("N137").Value = -12013.5

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20661.086
$ws.Range("I32").Value = 19549.268
$ws.Range("J32").Value = 28555
$ws.Range("K32").Value = 19549.268
$ws.Range("L32").Value = 28555
$ws.Range("M32").Value = -19262.268
$ws.Range("N32").Value = -29129

$ws.Range("H88").Value = 2718.8333
$ws.Range("I88").Value = 2489.2666
$ws.Range("J88").Value = 3866.6667
$ws.Range("K88").Value = 2489.2666
$ws.Range("L88").Value = 3866.6667
$ws.Range("M88").Value = -2083.2666
$ws.Range("N88").Value = -4678.6667

$ws.Range("H91").Value = 2718.8333
$ws.Range("I91").Value = 2489.2666
$ws.Range("J91").Value = 3866.6667
$ws.Range("K91").Value = 2489.2666
$ws.Range("L91").Value = 3866.6667
$ws.Range("M91").Value = -1085.2666
$ws.Range("N91").Value = -6674.6667

$ws.Range("H132").Value = 1865.8298
$ws.Range("I132").Value = 1623.72
$ws.Range("J132").Value = 2140.9546
$ws.Range("K132").Value = 4871.16
$ws.Range("L132").Value = 6422.8638
$ws.Range("M132").Value = -2341.16
$ws.Range("N132").Value = -11482.8638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 12108.083
$ws.Range("I107").Value = 1142.7142
$ws.Range("J107").Value = 27459.6
$ws.Range("K107").Value = 1142.7142
$ws.Range("L107").Value = 27459.6
$ws.Range("M107").Value = 777.2858000000001
$ws.Range("N107").Value = -31299.6

$ws.Range("H134").Value = 61328.383
$ws.Range("I134").Value = 2744.8076
$ws.Range("J134").Value = 251725
$ws.Range("K134").Value = 8234.4228
$ws.Range("L134").Value = 755175
$ws.Range("M134").Value = -5699.4228
$ws.Range("N134").Value = -760245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1618.871
$ws.Range("I31").Value = 834.2857
$ws.Range("J31").Value = 2265
$ws.Range("K31").Value = 834.2857
$ws.Range("L31").Value = 2265
$ws.Range("M31").Value = -539.2857
$ws.Range("N31").Value = -2855

$ws.Range("H34").Value = 1618.871
$ws.Range("I34").Value = 834.2857
$ws.Range("J34").Value = 2265
$ws.Range("K34").Value = 834.2857
$ws.Range("L34").Value = 2265
$ws.Range("M34").Value = -632.2857
$ws.Range("N34").Value = -2669

$ws.Range("H107").Value = 345.8421
$ws.Range("I107").Value = 265.23077
$ws.Range("J107").Value = 387.76
$ws.Range("K107").Value = 265.23077
$ws.Range("L107").Value = 387.76
$ws.Range("M107").Value = 1654.76923
$ws.Range("N107").Value = -4227.76

$ws.Range("H140").Value = 46747.145
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 46747.145
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 46747.145
$ws.Range("N140").Value = -57107.145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -23057

$ws.Range("H138").Value = 40642.223
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 40642.223
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 40642.223
$ws.Range("N138").Value = -50922.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 58000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 58000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 58000
$ws.Range("N135").Value = -68140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 1103
$ws.Range("I18").Value = 206
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 206
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -33
$ws.Range("N18").Value = -2346

$ws.Range("H24").Value = 4000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 4000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 4000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -4460

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H31").Value = 70012.664
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 70012.664
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 70012.664
$ws.Range("N31").Value = -70708.664

$ws.Range("H34").Value = 41014.5
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 41014.5
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 41014.5
$ws.Range("N34").Value = -41420.5

$ws.Range("H37").Value = 31676.334
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 31676.334
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 31676.334
$ws.Range("N37").Value = -32082.334

$ws.Range("H51").Value = 18000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 18000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 18000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -19020

$ws.Range("H55").Value = 12166.667
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 12166.667
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 12166.667
$ws.Range("N55").Value = -12720.667

$ws.Range("H123").Value = 34777.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 34777.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 34777.5
$ws.Range("N123").Value = -44577.5

$ws.Range("H136").Value = 2273.3035
$ws.Range("I136").Value = 2624.3103
$ws.Range("J136").Value = 1896.2963
$ws.Range("K136").Value = 7872.9309
$ws.Range("L136").Value = 5688.8889
$ws.Range("M136").Value = -5322.9309
$ws.Range("N136").Value = -10788.8889
